$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.598.16"
$ws.Range("E2").Value = "  +1.83%  "
$ws.Range("D3").Value = "1.829.00"
$ws.Range("E3").Value = "  +1.89%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'317.65"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "'0.5414"
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("D8").Value = "'0.4016"
$ws.Range("E8").Value = "  +6.28%  "
$ws.Range("D9").Value = "'0.07716"
$ws.Range("E9").Value = "  +3.57%  "
$ws.Range("D10").Value = "'1.123"
$ws.Range("E10").Value = "  +2.64%  "
$ws.Range("D11").Value = "'41.91"
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("D12").Value = "'21.36"
$ws.Range("E12").Value = "  +3.99%  "
$ws.Range("D13").Value = "'6.341"
$ws.Range("E13").Value = "  +3.73%  "
$ws.Range("D14").Value = "'7.645"
$ws.Range("E14").Value = "  +5.65%  "
$ws.Range("D15").Value = "'1.000"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").Value = "1.827.44"
$ws.Range("E16").Value = "  +1.83%  "
$ws.Range("D17").Value = "'0.00001091"
$ws.Range("E17").Value = "  +2.93%  "
$ws.Range("D18").Value = "'90.11"
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("D19").Value = "'0.06594"
$ws.Range("E19").Value = "  +1.59%  "
$ws.Range("D20").Value = "'17.84"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").Value = "'6.080"
$ws.Range("E22").Value = "  +3.10%  "
$ws.Range("D23").Value = "28.605.44"
$ws.Range("E23").Value = "  +1.76%  "
$ws.Range("E24").Value = "  +0.29%  "
$ws.Range("D25").Value = "'2.274"
$ws.Range("E25").Value = "  +8.82%  "
$ws.Range("D26").Value = "'158.40"
$ws.Range("E26").Value = "  +1.96%  "
$ws.Range("D27").Value = "'2.464"
$ws.Range("E27").Value = "  +7.93%  "
$ws.Range("D28").Value = "'20.81"
$ws.Range("E28").Value = "  +2.55%  "
$ws.Range("D29").Value = "2.038.55"
$ws.Range("E29").Value = "  +2.00%  "
$ws.Range("D30").Value = "'124.38"
$ws.Range("E30").Value = "  +2.67%  "
$ws.Range("D31").Value = "'1.130"
$ws.Range("E31").Value = "  +1.22%  "
$ws.Range("D32").Value = "'0.1113"
$ws.Range("E32").Value = "  +4.68%  "
$ws.Range("D33").Value = "'0.07589"
$ws.Range("E33").Value = "  +17.33%  "
$ws.Range("D34").Value = "'5.699"
$ws.Range("E34").Value = "  +2.70%  "
$ws.Range("D35").Value = "'3.648"
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("D36").Value = "'0.2254"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").Value = "'0.02361"
$ws.Range("E37").Value = "  +3.12%  "
$ws.Range("D38").Value = "'8.942"
$ws.Range("E38").Value = "  +5.93%  "
$ws.Range("D39").Value = "'5.220"
$ws.Range("E39").Value = "  +4.22%  "
$ws.Range("D40").Value = "'0.6312"
$ws.Range("E40").Value = "  +2.16%  "
$ws.Range("E41").Value = "  +2.46%  "
$ws.Range("D42").Value = "'1.190"
$ws.Range("E42").Value = "  +1.08%  "
$ws.Range("D43").Value = "'0.9997"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").Value = "'1.401"
$ws.Range("E44").Value = "  -3.43%  "
$ws.Range("D45").Value = "'13.53"
$ws.Range("E45").Value = "  +2.03%  "
$ws.Range("D46").Value = "'0.5904"
$ws.Range("E46").Value = "  +2.20%  "
$ws.Range("D47").Value = "'3.711"
$ws.Range("E47").Value = "  +1.04%  "
$ws.Range("D48").Value = "'125.52"
$ws.Range("E48").Value = "  +1.03%  "
$ws.Range("D49").Value = "'2.006"
$ws.Range("E49").Value = "  +4.25%  "
$ws.Range("E50").Value = "  +0.90%  "
$ws.Range("D51").Value = "'0.06917"
$ws.Range("E51").Value = "  +1.49%  "
